# Updates the cryptos list (Price / Volume(1h) columns) to the latest
# scraped snapshot. Price values that look like plain numbers are written
# with a leading apostrophe so Excel keeps them as literal text (matching
# the sheet's original t="inlineStr" cells) instead of silently parsing
# them into numbers - which would also truncate significant trailing
# zeros (e.g. "0.0840" -> 0.084). ClearFormats() then drops the "number
# stored as text" quote-prefix flag Excel auto-applies, so the cell style
# stays the same untouched default as every other data cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.770.54'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '2.099.96'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = "'227.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").Value = "'61.89"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.83%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = "'0.388"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").Value = "'0.0840"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").Value = "'15.76"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.79%  '
$ws.Range("D13").Value = '2.412.91'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").Value = "'0.799"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = "'5.47"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '2.107.32'
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").Value = '38.759.89'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = "'71.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = "'6.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '0.0₃0842'
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").Value = "'226.92"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = "'2.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.23%  '
$ws.Range("D25").Value = "'2.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("D26").Value = "'9.65"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("D27").Value = "'170.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").Value = "'0.136"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("D29").Value = "'1.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.44%  '
$ws.Range("D30").Value = "'19.33"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("E31").Value = '  +8.85%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").Value = "'4.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("D34").Value = "'4.79"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").Value = "'7.11"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +11.04%  '
$ws.Range("D36").Value = "'0.0612"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = "'2.35"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.96%  '
$ws.Range("D38").Value = "'3.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").Value = "'17.98"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.91%  '
$ws.Range("E41").Value = '  +2.64%  '
$ws.Range("D42").Value = "'101.52"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").Value = '1.525.03'
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("E44").Value = '  +7.81%  '
$ws.Range("D45").Value = "'2.80"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("D46").Value = "'7.76"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("D47").Value = "'0.0910"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.94%  '
$ws.Range("E48").Value = '  +4.42%  '
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("E50").Value = '  -1.38%  '
$ws.Range("D51").Value = '2.300.16'
$ws.Range("E51").Value = '  +0.40%  '
